# feat: add 2022-Q1 data
#
# The workbook currently has sheets:
#   2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
#
# Target layout:
#   2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# The existing "总计" sheet becomes "2022-Q1" (reusing its sheetId/rId) and is
# repopulated with the quarterly fund-holding detail table. A brand-new "总计"
# sheet is appended at the end, holding the historical summary table plus a
# new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$fundSheet = $wb.Worksheets.Item("2021-Q4")
$oldTotal  = $wb.Worksheets.Item("总计")

# Grab formatting references before anything is overwritten.
$oldTotal.Range("B1:D1").Copy()
$totalHeaderFmtSrc = $true

# ---------------------------------------------------------------------
# 1) Create the new "总计" sheet (appended after the current total sheet)
#    and fill it with the historical summary plus the 2022-Q1 row.
# ---------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add($null, $oldTotal)
$newTotal.Name = "总计__tmp"

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"
$newTotal.Range("B1:D1").PasteSpecial(-4122)

$totalRows = @(
    @(0, "2022-Q1", 15, 5.11),
    @(1, "2021-Q4", 26, 5.58),
    @(2, "2021-Q3", 21, 5),
    @(3, "2021-Q2", 41, 11.28),
    @(4, "2021-Q1", 75, 16.86),
    @(5, "2020-Q4", 49, 12.53)
)

$r = 2
foreach ($row in $totalRows) {
    $newTotal.Cells.Item($r, 1).Value = $row[0]
    $newTotal.Cells.Item($r, 2).Value = $row[1]
    $newTotal.Cells.Item($r, 3).Value = $row[2]
    $newTotal.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$oldTotal.Range("A2").Copy()
$newTotal.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Turn the old "总计" sheet into "2022-Q1" and replace its contents
#    with the fund-holding detail table.
# ---------------------------------------------------------------------
$oldTotal.Range("A1:D6").ClearContents()
$oldTotal.Name = "2022-Q1"
$q1 = $oldTotal

$fundSheet.Range("B1:H1").Copy()
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$q1.Range("B1:H1").PasteSpecial(-4122)

$fundRows = @(
    @(0, "000991", "工银瑞信战略转型主题股票A", "48.06", "88.32", "4.94", "2.3742", 4),
    @(1, "011521", "鹏扬景源一年持有期混合A",   "33.78", "21.39", "2.26", "0.7634", 1),
    @(2, "000971", "诺安新经济股票",             "15.43", "82.95", "4.46", "0.6882", 7),
    @(3, "009064", "鹏扬景沃六个月持有期混合A", "23.92", "21.80", "1.93", "0.4617", 2),
    @(4, "011473", "工银瑞信战略转型主题股票C", "5.61",  "88.32", "4.94", "0.2771", 4),
    @(5, "009065", "鹏扬景沃六个月持有期混合C", "7.70",  "21.80", "1.93", "0.1486", 2),
    @(6, "393001", "中海优势精选灵活配置混合",   "1.58",  "78.65", "7.68", "0.1213", 6),
    @(7, "011522", "鹏扬景源一年持有期混合C",   "4.58",  "21.39", "2.26", "0.1035", 1),
    @(8, "005576", "华泰柏瑞新金融地产灵活配置混合", "0.79", "94.50", "6.59", "0.0521", 5),
    @(9, "009927", "工银瑞信聚利18个月定期开放混合A", "5.54", "23.27", "0.79", "0.0438", 10),
    @(10, "001780", "诺安改革趋势灵活配置混合",  "0.46",  "68.34", "8.48", "0.0390", 2),
    @(11, "510060", "工银上证央企50ETF",          "0.80",  "99.17", "2.80", "0.0224", 9),
    @(12, "009928", "工银瑞信聚利18个月定期开放混合C", "0.83", "23.27", "0.79", "0.0066", 10),
    @(13, "004988", "人保双利优选混合A",          "0.58",  "25.37", "0.57", "0.0033", 7),
    @(14, "004989", "人保双利优选混合C",          "0.00",  "25.37", "0.57", $null, 7)
)

$q1.Range("B2:H16").NumberFormat = "@"

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    if ($row[6] -eq $null) {
        $q1.Cells.Item($r, 7).NumberFormat = "General"
        $q1.Cells.Item($r, 7).Value = 0
    } else {
        $q1.Cells.Item($r, 7).Value = $row[6]
    }
    $q1.Cells.Item($r, 8).NumberFormat = "General"
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$fundSheet.Range("A2:A16").Copy()
$q1.Range("A2:A16").PasteSpecial(-4122)

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Final rename of the appended sheet to "总计".
# ---------------------------------------------------------------------
$newTotal.Name = "总计"
